$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.01609314246516
$ws.Range("D2").Value = 1.04613136794257
$ws.Range("E2").Value = 1.017648189476872
$ws.Range("F2").Value = 1.048039328528655
$ws.Range("I2").Value = 1.035654397595778
$ws.Range("J2").Value = 1.021314919512776
$ws.Range("K2").Value = 1.048897663452835
$ws.Range("L2").Value = 1.020496754285191
$ws.Range("M2").Value = 1.050800286572524
$ws.Range("N2").Value = 1.022765303844943

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.017204715642761
$ws.Range("D3").Value = 1.046789464402056
$ws.Range("E3").Value = 1.018596214336293
$ws.Range("F3").Value = 1.048948339898391
$ws.Range("I3").Value = 1.035763180007724
$ws.Range("J3").Value = 1.022060939589109
$ws.Range("K3").Value = 1.049367739711606
$ws.Range("L3").Value = 1.021250026443115
$ws.Range("M3").Value = 1.051521011763711
$ws.Range("N3").Value = 1.023512383355354

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.017924426576686
$ws.Range("D4").Value = 1.047214778117057
$ws.Range("E4").Value = 1.01921038460264
$ws.Range("F4").Value = 1.049536499637218
$ws.Range("I4").Value = 1.035832007948652
$ws.Range("J4").Value = 1.022543613419927
$ws.Range("K4").Value = 1.049670708746801
$ws.Range("L4").Value = 1.021737563420935
$ws.Range("M4").Value = 1.051986707106567
$ws.Range("N4").Value = 1.023995742638375

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.018227101033186
$ws.Range("D5").Value = 1.04739345412321
$ws.Range("E5").Value = 1.019468757428537
$ws.Range("F5").Value = 1.049783753026784
$ws.Range("I5").Value = 1.03586056872539
$ws.Range("J5").Value = 1.022746517219511
$ws.Range("K5").Value = 1.049797787896901
$ws.Range("L5").Value = 1.021942552434032
$ws.Range("M5").Value = 1.052182325619244
$ws.Range("N5").Value = 1.024198934584628

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.018277927754283
$ws.Range("D6").Value = 1.047423447179581
$ws.Range("E6").Value = 1.019512149655244
$ws.Range("F6").Value = 1.049825267391012
$ws.Range("I6").Value = 1.035865342226266
$ws.Range("J6").Value = 1.02278058493859
$ws.Range("K6").Value = 1.049819108057481
$ws.Range("L6").Value = 1.021976972656489
$ws.Range("M6").Value = 1.052215161418703
$ws.Range("N6").Value = 1.024233050683776

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.017928470502606
$ws.Range("D7").Value = 1.047217166091386
$ws.Range("E7").Value = 1.019213836303148
$ws.Range("F7").Value = 1.04953980348618
$ws.Range("I7").Value = 1.035832391051462
$ws.Range("J7").Value = 1.02254632467875
$ws.Range("K7").Value = 1.049672407922027
$ws.Range("L7").Value = 1.021740302382486
$ws.Range("M7").Value = 1.051989321600327
$ws.Range("N7").Value = 1.023998457747497

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.016468711233159
$ws.Range("D8").Value = 1.046353881101592
$ws.Range("E8").Value = 1.017968426614403
$ws.Range("F8").Value = 1.04834653881018
$ws.Range("I8").Value = 1.035691484060837
$ws.Range("J8").Value = 1.02156705102337
$ws.Range("K8").Value = 1.049056775952822
$ws.Range("L8").Value = 1.020751301301716
$ws.Range("M8").Value = 1.051043995318694
$ws.Range("N8").Value = 1.023017793411202

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.013899841300623
$ws.Range("D9").Value = 1.04482876483267
$ws.Range("E9").Value = 1.015779494598929
$ws.Range("F9").Value = 1.046243677052117
$ws.Range("I9").Value = 1.035431256894698
$ws.Range("J9").Value = 1.019841057100759
$ws.Range("K9").Value = 1.047962808488495
$ws.Range("L9").Value = 1.019009476606122
$ws.Range("M9").Value = 1.049373196786523
$ws.Range("N9").Value = 1.021289348379247

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.012189518994097
$ws.Range("D10").Value = 1.04380949996781
$ws.Range("E10").Value = 1.014324010994939
$ws.Range("F10").Value = 1.044841734426725
$ws.Range("I10").Value = 1.035249788504611
$ws.Range("J10").Value = 1.018690134154438
$ws.Range("K10").Value = 1.047227426216062
$ws.Range("L10").Value = 1.017848886053107
$ws.Range("M10").Value = 1.048256037142623
$ws.Range("N10").Value = 1.020136790990329

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.01144945713167
$ws.Range("D11").Value = 1.043367569972488
$ws.Range("E11").Value = 1.01369467490555
$ws.Range("F11").Value = 1.044234686810687
$ws.Range("I11").Value = 1.035169324985274
$ws.Range("J11").Value = 1.018191709556401
$ws.Range("K11").Value = 1.04690757390006
$ws.Range("L11").Value = 1.017346486650415
$ws.Range("M11").Value = 1.047771529075733
$ws.Range("N11").Value = 1.019637658572192

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.011174642316126
$ws.Range("D12").Value = 1.043203331763366
$ws.Range("E12").Value = 1.013461046442278
$ws.Range("F12").Value = 1.044009204051422
$ws.Range("I12").Value = 1.035139154429374
$ws.Range("J12").Value = 1.0180065622172
$ws.Range("K12").Value = 1.046788553389632
$ws.Range("L12").Value = 1.017159894743835
$ws.Range("M12").Value = 1.047591446643303
$ws.Range("N12").Value = 1.019452248302533

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.011233587551235
$ws.Range("D13").Value = 1.043238565269972
$ws.Range("E13").Value = 1.013511154431224
$ws.Range("F13").Value = 1.044057570767941
$ws.Range("I13").Value = 1.035145638904648
$ws.Range("J13").Value = 1.018046277415904
$ws.Range("K13").Value = 1.046814093317007
$ws.Range("L13").Value = 1.017199918364145
$ws.Range("M13").Value = 1.047630080101298
$ws.Range("N13").Value = 1.019492019901375

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01142673928519
$ws.Range("D14").Value = 1.043353995730312
$ws.Range("E14").Value = 1.013675360349845
$ws.Range("F14").Value = 1.044216048291981
$ws.Range("I14").Value = 1.03516683684371
$ws.Range("J14").Value = 1.018176405429448
$ws.Range("K14").Value = 1.046897739967819
$ws.Range("L14").Value = 1.017331062464714
$ws.Range("M14").Value = 1.047756645735799
$ws.Range("N14").Value = 1.019622332711624

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.011545756635687
$ws.Range("D15").Value = 1.043425104935347
$ws.Range("E15").Value = 1.013776550924567
$ws.Range("F15").Value = 1.044313691797242
$ws.Range("I15").Value = 1.035179860134948
$ws.Range("J15").Value = 1.018256580231138
$ws.Range("K15").Value = 1.046949249227421
$ws.Range("L15").Value = 1.017411867534896
$ws.Range("M15").Value = 1.047834611855593
$ws.Range("N15").Value = 1.019702621370729

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.01223864532341
$ws.Range("D16").Value = 1.043838817293395
$ws.Range("E16").Value = 1.014365796883352
$ws.Range("F16").Value = 1.044882022350592
$ws.Range("I16").Value = 1.035255088908458
$ws.Range("J16").Value = 1.018723211522991
$ws.Range("K16").Value = 1.047248623804812
$ws.Range("L16").Value = 1.017882231677582
$ws.Range("M16").Value = 1.048288176215346
$ws.Range("N16").Value = 1.020169915332539

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.012673414429248
$ws.Range("D17").Value = 1.044098173582615
$ws.Range("E17").Value = 1.014735655759211
$ws.Range("F17").Value = 1.04523852278291
$ws.Range("I17").Value = 1.035301773285984
$ws.Range("J17").Value = 1.019015898875742
$ws.Range("K17").Value = 1.047436032429017
$ws.Range("L17").Value = 1.018177317200724
$ws.Range("M17").Value = 1.048572479704704
$ws.Range("N17").Value = 1.020463018334902

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.012927058024637
$ws.Range("D18").Value = 1.044249395493985
$ws.Range("E18").Value = 1.014951474814405
$ws.Range("F18").Value = 1.045446463633679
$ws.Range("I18").Value = 1.035328821432371
$ws.Range("J18").Value = 1.019186611965337
$ws.Range("K18").Value = 1.047545206945481
$ws.Range("L18").Value = 1.018349449486294
$ws.Range("M18").Value = 1.048738234690639
$ws.Range("N18").Value = 1.020633973856666

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.013013552465369
$ws.Range("D19").Value = 1.044300948671112
$ws.Range("E19").Value = 1.01502507822658
$ws.Range("F19").Value = 1.045517366080009
$ws.Range("I19").Value = 1.035338013252149
$ws.Range("J19").Value = 1.019244819620204
$ws.Range("K19").Value = 1.047582409228731
$ws.Range("L19").Value = 1.018408144504654
$ws.Range("M19").Value = 1.048794740196104
$ws.Range("N19").Value = 1.02069226417308

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.012626762639507
$ws.Range("D20").Value = 1.044070352910907
$ws.Range("E20").Value = 1.014695964440876
$ws.Range("F20").Value = 1.04520027363206
$ws.Range("I20").Value = 1.035296783321148
$ws.Range("J20").Value = 1.018984496965033
$ws.Range("K20").Value = 1.047415939497595
$ws.Range("L20").Value = 1.018145655879913
$ws.Range("M20").Value = 1.048541984320715
$ws.Range("N20").Value = 1.020431571829879

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.011369858791944
$ws.Range("D21").Value = 1.043320006690225
$ws.Range("E21").Value = 1.013627002072736
$ws.Range("F21").Value = 1.044169380551981
$ws.Range("I21").Value = 1.035160602381434
$ws.Range("J21").Value = 1.018138086256828
$ws.Range("K21").Value = 1.046873113987005
$ws.Range("L21").Value = 1.017292443199704
$ws.Range("M21").Value = 1.047719378463863
$ws.Range("N21").Value = 1.019583959121383

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.010580039178075
$ws.Range("D22").Value = 1.04284773902988
$ws.Range("E22").Value = 1.012955684007076
$ws.Range("F22").Value = 1.043521227962742
$ws.Range("I22").Value = 1.035073343995386
$ws.Range("J22").Value = 1.017605854878073
$ws.Range("K22").Value = 1.046530585624628
$ws.Range("L22").Value = 1.016756119487359
$ws.Range("M22").Value = 1.047201511275958
$ws.Range("N22").Value = 1.019050971913023

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.010998695535732
$ws.Range("D23").Value = 1.043098143445499
$ws.Range("E23").Value = 1.013311488358428
$ws.Range("F23").Value = 1.043864824445042
$ws.Range("I23").Value = 1.03511975619175
$ws.Range("J23").Value = 1.017888006472803
$ws.Range("K23").Value = 1.046712282789473
$ws.Range("L23").Value = 1.017040422998756
$ws.Range("M23").Value = 1.047476104841717
$ws.Range("N23").Value = 1.019333524195381

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.012647842419988
$ws.Range("D24").Value = 1.044082924048791
$ws.Range("E24").Value = 1.014713898971354
$ws.Range("F24").Value = 1.04521755677642
$ws.Range("I24").Value = 1.03529903863405
$ws.Range("J24").Value = 1.018998686156719
$ws.Range("K24").Value = 1.04742501905417
$ws.Range("L24").Value = 1.018159962224898
$ws.Range("M24").Value = 1.048555764102107
$ws.Range("N24").Value = 1.020445781171845

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.014563555161645
$ws.Range("D25").Value = 1.045223495524088
$ws.Range("E25").Value = 1.016344716402859
$ws.Range("F25").Value = 1.046787329498541
$ws.Range("I25").Value = 1.035499941805617
$ws.Range("J25").Value = 1.020287313863202
$ws.Range("K25").Value = 1.048246700988321
$ws.Range("L25").Value = 1.019459670186604
$ws.Range("M25").Value = 1.049805724413429
$ws.Range("N25").Value = 1.021736238877479
